$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy formatting (borders/fill/wrap) from row 11 down to rows 12-15 (new rows)
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E15").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = "TestCase_E1"
$ws.Range("B2").Value = "OPQA-256"
$ws.Range("C2").Value = "Verify that user is able to add an Article from ALL content search results page to a particular watchlist"
$prefixLen = "Verify that user is able to add an Article from ".Length
$boldLen = "ALL".Length
$bChars = $ws.Range("C2").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add an Article from ALL content search results page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C2").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D2").Value = "N"
$ws.Range("E2").Value = "SKIP"

# Row 3
$ws.Range("A3").Value = "TestCase_E2"
$ws.Range("B3").Value = "OPQA-259"
$ws.Range("C3").Value = "Verify that user is able to add a Patent from ALL content search results page to a particular watchlist"
$prefixLen = "Verify that user is able to add a Patent from ".Length
$boldLen = "ALL".Length
$bChars = $ws.Range("C3").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add a Patent from ALL content search results page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C3").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D3").Value = "N"
$ws.Range("E3").Value = "SKIP"

# Row 4
$ws.Range("A4").Value = "TestCase_E3"
$ws.Range("B4").Value = "OPQA-260"
$ws.Range("C4").Value = "Verify that user is able to add a Post from ALL content search results page to a particular watchlist"
$prefixLen = "Verify that user is able to add a Post from ".Length
$boldLen = "ALL".Length
$bChars = $ws.Range("C4").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add a Post from ALL content search results page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C4").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D4").Value = "N"
$ws.Range("E4").Value = "SKIP"

# Row 5
$ws.Range("A5").Value = "TestCase_E4"
$ws.Range("B5").Value = "OPQA-261"
$ws.Range("C5").Value = "Verify that user is able to unwatch an Article from ALL content search results page"
$prefixLen = "Verify that user is able to unwatch an Article from ".Length
$boldLen = "ALL".Length
$bChars = $ws.Range("C5").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch an Article from ALL content search results page".Length - $restStart + 1
$restChars = $ws.Range("C5").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D5").Value = "N"
$ws.Range("E5").Value = "SKIP"

# Row 6
$ws.Range("A6").Value = "TestCase_E5"
$ws.Range("B6").Value = "OPQA-262"
$ws.Range("C6").Value = "Verify that user is able to unwatch a Patent from ALL content search results page"
$prefixLen = "Verify that user is able to unwatch a Patent from ".Length
$boldLen = "ALL".Length
$bChars = $ws.Range("C6").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch a Patent from ALL content search results page".Length - $restStart + 1
$restChars = $ws.Range("C6").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D6").Value = "N"
$ws.Range("E6").Value = "SKIP"

# Row 7
$ws.Range("A7").Value = "TestCase_E6"
$ws.Range("B7").Value = "OPQA-264"
$ws.Range("C7").Value = "Verify that user is able to unwatch a Post from ALL content search results page"
$prefixLen = "Verify that user is able to unwatch a Post from ".Length
$boldLen = "ALL".Length
$bChars = $ws.Range("C7").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch a Post from ALL content search results page".Length - $restStart + 1
$restChars = $ws.Range("C7").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D7").Value = "N"
$ws.Range("E7").Value = "SKIP"

# Row 8
$ws.Range("A8").Value = "TestCase_E7"
$ws.Range("B8").Value = "OPQA-265"
$ws.Range("C8").Value = "Verify that user is able to add an Article from Articles content search results page to a particular watchlist"
$prefixLen = "Verify that user is able to add an Article from ".Length
$boldLen = "Articles".Length
$bChars = $ws.Range("C8").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add an Article from Articles content search results page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C8").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D8").Value = "N"
$ws.Range("E8").Value = "SKIP"

# Row 9
$ws.Range("A9").Value = "TestCase_E8"
$ws.Range("B9").Value = "OPQA-267"
$ws.Range("C9").Value = "Verify that user is able to unwatch an Article from Articles content search results page"
$prefixLen = "Verify that user is able to unwatch an Article from ".Length
$boldLen = "Articles".Length
$bChars = $ws.Range("C9").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch an Article from Articles content search results page".Length - $restStart + 1
$restChars = $ws.Range("C9").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D9").Value = "N"
$ws.Range("E9").Value = "SKIP"

# Row 10
$ws.Range("A10").Value = "TestCase_E9"
$ws.Range("B10").Value = "OPQA-268"
$ws.Range("C10").Value = "Verify that user is able to add an Article from Record View page to a particular watchlist"
$prefixLen = "Verify that user is able to add an Article from ".Length
$boldLen = "Record View ".Length
$bChars = $ws.Range("C10").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add an Article from Record View page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C10").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D10").Value = "N"
$ws.Range("E10").Value = "SKIP"

# Row 11
$ws.Range("A11").Value = "TestCase_E10"
$ws.Range("B11").Value = "OPQA-269"
$ws.Range("C11").Value = "Verify that user is able to unwatch an Article from Record View page"
$prefixLen = "Verify that user is able to unwatch an Article from ".Length
$boldLen = "Record View".Length
$bChars = $ws.Range("C11").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch an Article from Record View page".Length - $restStart + 1
$restChars = $ws.Range("C11").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D11").Value = "N"
$ws.Range("E11").Value = "SKIP"

# Row 12
$ws.Range("A12").Value = "TestCase_E11"
$ws.Range("B12").Value = "OPQA-272"
$ws.Range("C12").Value = "Verify that user is able to add a Patent from Patents content search results page to a particular watchlist"
$prefixLen = "Verify that user is able to add a Patent from ".Length
$boldLen = "Patents".Length
$bChars = $ws.Range("C12").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add a Patent from Patents content search results page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C12").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D12").Value = "N"
$ws.Range("E12").Value = "SKIP"

# Row 13
$ws.Range("A13").Value = "TestCase_E12"
$ws.Range("B13").Value = "OPQA-273"
$ws.Range("C13").Value = "Verify that user is able to unwatch a Patent from Patents content search results page"
$prefixLen = "Verify that user is able to unwatch a Patent from ".Length
$boldLen = "Patents".Length
$bChars = $ws.Range("C13").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch a Patent from Patents content search results page".Length - $restStart + 1
$restChars = $ws.Range("C13").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D13").Value = "N"
$ws.Range("E13").Value = "SKIP"

# Row 14
$ws.Range("A14").Value = "TestCase_E13"
$ws.Range("B14").Value = "OPQA-276"
$ws.Range("C14").Value = "Verify that user is able to add a Patent from Record View page to a particular watchlist"
$prefixLen = "Verify that user is able to add a Patent from ".Length
$boldLen = "Record View ".Length
$bChars = $ws.Range("C14").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to add a Patent from Record View page to a particular watchlist".Length - $restStart + 1
$restChars = $ws.Range("C14").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D14").Value = "N"
$ws.Range("E14").Value = "SKIP"

# Row 15
$ws.Range("A15").Value = "TestCase_E14"
$ws.Range("B15").Value = "OPQA-277"
$ws.Range("C15").Value = "Verify that user is able to unwatch a Patent from Record View page"
$prefixLen = "Verify that user is able to unwatch a Patent from ".Length
$boldLen = "Record View".Length
$bChars = $ws.Range("C15").Characters($prefixLen + 1, $boldLen)
$bChars.Font.Bold = $true
$bChars.Font.ColorIndex = 8
$restStart = $prefixLen + $boldLen + 1
$restLen = "Verify that user is able to unwatch a Patent from Record View page".Length - $restStart + 1
$restChars = $ws.Range("C15").Characters($restStart, $restLen)
$restChars.Font.Size = 11
$ws.Range("D15").Value = "Y"
$ws.Range("E15").Value = "PASS"

$ws.Range("C15").Select()